$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 8, shifting rows 8:38 down to 9:39
$ws.Rows.Item(8).Insert()

# Fill in the new row 8 with the same "template" values as the other
# Guayaba / Vega Modelo de Temuco rows, but with the new date/volume/price data.
$ws.Cells.Item(8, 1).Value = 10
$ws.Cells.Item(8, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(8, 3).Value = "La Araucanía"
$ws.Cells.Item(8, 4).Value = 45099
$ws.Cells.Item(8, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(8, 5).Value = 9
$ws.Cells.Item(8, 6).Value = "Fruta"
$ws.Cells.Item(8, 7).Value = 100108
$ws.Cells.Item(8, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(8, 9).Value = 100108001
$ws.Cells.Item(8, 10).Value = "Guayaba"
$ws.Cells.Item(8, 11).Value = "Sin especificar"
$ws.Cells.Item(8, 12).Value = "Primera"
$ws.Cells.Item(8, 13).Value = 200
$ws.Cells.Item(8, 14).Value = 2600
$ws.Cells.Item(8, 15).Value = 2600
$ws.Cells.Item(8, 16).Value = 2600
$ws.Cells.Item(8, 17).Value = "`$/kilo"
$ws.Cells.Item(8, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(8, 19).Value = 2600
$ws.Cells.Item(8, 20).Value = 1
